$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / Modified date) from 45183 to 45184 for rows 2-70
for ($row = 2; $row -le 70; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
